$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Raw")

$ws.Range("A39").Value = 'abt'
$ws.Range("B39").Value = '2021-01-11 01:33:27.922'
$ws.Range("A40").Value = 'alk'
$ws.Range("B40").Value = '2021-01-11 01:33:28.031'
$ws.Range("A41").Value = 'alb'
$ws.Range("B41").Value = '2021-01-11 01:33:28.056'
$ws.Range("A42").Value = 'all'
$ws.Range("B42").Value = '2021-01-11 01:33:28.079'
$ws.Range("A43").Value = 'amgn'
$ws.Range("B43").Value = '2021-01-11 01:33:28.1'
$ws.Range("A44").Value = 'anh'
$ws.Range("B44").Value = '2021-01-11 01:33:28.12'
$ws.Range("A45").Value = 'aapl'
$ws.Range("B45").Value = '2021-01-11 01:33:28.141'
$ws.Range("A46").Value = 'adm'
$ws.Range("B46").Value = '2021-01-11 01:33:28.161'
$ws.Range("A47").Value = 'aca'
$ws.Range("B47").Value = '2021-01-11 01:33:28.183'
$ws.Range("A48").Value = 't'
$ws.Range("B48").Value = '2021-01-11 01:33:28.201'
$ws.Range("A49").Value = 'bac'
$ws.Range("B49").Value = '2021-01-11 01:33:28.222'
$ws.Range("A50").Value = 'bk'
$ws.Range("B50").Value = '2021-01-11 01:33:28.242'
$ws.Range("A51").Value = 'bhe'
$ws.Range("B51").Value = '2021-01-11 01:33:28.265'
$ws.Range("A52").Value = 'bhp'
$ws.Range("B52").Value = '2021-01-11 01:33:28.283'
$ws.Range("A53").Value = 'biib'
$ws.Range("B53").Value = '2021-01-11 01:33:28.304'
$ws.Range("A54").Value = 'bhf'
$ws.Range("B54").Value = '2021-01-11 01:33:28.322'
$ws.Range("A55").Value = 'cm'
$ws.Range("B55").Value = '2021-01-11 01:33:28.339'
$ws.Range("A56").Value = 'cof'
$ws.Range("B56").Value = '2021-01-11 01:33:28.356'
$ws.Range("A57").Value = 'cat'
$ws.Range("B57").Value = '2021-01-11 01:33:28.372'
$ws.Range("A58").Value = 'ce'
$ws.Range("B58").Value = '2021-01-11 01:33:28.388'
$ws.Range("A59").Value = 'csco'
$ws.Range("B59").Value = '2021-01-11 01:33:28.408'
$ws.Range("A60").Value = 'cohu'
$ws.Range("B60").Value = '2021-01-11 01:33:28.424'
$ws.Range("A61").Value = 'cmcsa'
$ws.Range("B61").Value = '2021-01-11 01:33:28.445'
$ws.Range("A62").Value = 'glw'
$ws.Range("B62").Value = '2021-01-11 01:33:28.463'
$ws.Range("A63").Value = 'cmi'
$ws.Range("B63").Value = '2021-01-11 01:33:28.478'
$ws.Range("A64").Value = 'cvs'
$ws.Range("B64").Value = '2021-01-11 01:33:28.492'
$ws.Range("A65").Value = 'de'
$ws.Range("B65").Value = '2021-01-11 01:33:28.506'
$ws.Range("A66").Value = 'dal'
$ws.Range("B66").Value = '2021-01-11 01:33:28.52'
$ws.Range("A67").Value = 'dlr'
$ws.Range("B67").Value = '2021-01-11 01:33:28.535'
$ws.Range("A68").Value = 'etn'
$ws.Range("B68").Value = '2021-01-11 01:33:28.549'
$ws.Range("A69").Value = 'xom'
$ws.Range("B69").Value = '2021-01-11 01:33:28.562'
$ws.Range("A70").Value = 'fdx'
$ws.Range("B70").Value = '2021-01-11 01:33:28.578'
$ws.Range("A71").Value = 'fitb'
$ws.Range("B71").Value = '2021-01-11 01:33:28.592'
$ws.Range("A72").Value = 'flr'
$ws.Range("B72").Value = '2021-01-11 01:33:28.607'
$ws.Range("A73").Value = 'fl'
$ws.Range("B73").Value = '2021-01-11 01:33:28.623'
$ws.Range("A74").Value = 'gm'
$ws.Range("B74").Value = '2021-01-11 01:33:28.639'
$ws.Range("A75").Value = 'gild'
$ws.Range("B75").Value = '2021-01-11 01:33:28.653'
$ws.Range("A76").Value = 'gs'
$ws.Range("B76").Value = '2021-01-11 01:33:28.666'
$ws.Range("A77").Value = 'gt'
$ws.Range("B77").Value = '2021-01-11 01:33:28.678'
$ws.Range("A78").Value = 'hal'
$ws.Range("B78").Value = '2021-01-11 01:33:28.7'
$ws.Range("A79").Value = 'hfc'
$ws.Range("B79").Value = '2021-01-11 01:33:28.713'
$ws.Range("A80").Value = 'hmc'
$ws.Range("B80").Value = '2021-01-11 01:33:28.724'
$ws.Range("A81").Value = 'hsbc'
$ws.Range("B81").Value = '2021-01-11 01:33:28.735'
$ws.Range("A82").Value = 'intc'
$ws.Range("B82").Value = '2021-01-11 01:33:28.75'
$ws.Range("A83").Value = 'ibm'
$ws.Range("B83").Value = '2021-01-11 01:33:28.759'
$ws.Range("A84").Value = 'Ip'
$ws.Range("B84").Value = '2021-01-11 01:33:28.769'
$ws.Range("A85").Value = 'jbl'
$ws.Range("B85").Value = '2021-01-11 01:33:28.783'
$ws.Range("A86").Value = 'jnj'
$ws.Range("B86").Value = '2021-01-11 01:33:28.795'
$ws.Range("A87").Value = 'jpm'
$ws.Range("B87").Value = '2021-01-11 01:33:28.808'
$ws.Range("A88").Value = 'jnpr'
$ws.Range("B88").Value = '2021-01-11 01:33:28.82'
$ws.Range("A89").Value = 'key'
$ws.Range("B89").Value = '2021-01-11 01:33:28.832'
$ws.Range("A90").Value = 'kmb'
$ws.Range("B90").Value = '2021-01-11 01:33:28.844'
$ws.Range("A91").Value = 'kss'
$ws.Range("B91").Value = '2021-01-11 01:33:28.856'
$ws.Range("A92").Value = 'kr'
$ws.Range("B92").Value = '2021-01-11 01:33:28.87'
$ws.Range("A93").Value = 'man'
$ws.Range("B93").Value = '2021-01-11 01:33:28.883'
$ws.Range("A94").Value = 'mck'
$ws.Range("B94").Value = '2021-01-11 01:33:28.894'
$ws.Range("A95").Value = 'mdc'
$ws.Range("B95").Value = '2021-01-11 01:33:28.905'
$ws.Range("A96").Value = 'mdt'
$ws.Range("B96").Value = '2021-01-11 01:33:28.917'
$ws.Range("A97").Value = 'mrk'
$ws.Range("B97").Value = '2021-01-11 01:33:28.928'
$ws.Range("A98").Value = 'met'
$ws.Range("B98").Value = '2021-01-11 01:33:28.937'
$ws.Range("A99").Value = 'msft'
$ws.Range("B99").Value = '2021-01-11 01:33:28.947'
$ws.Range("A100").Value = 'mos'
$ws.Range("B100").Value = '2021-01-11 01:33:28.958'
$ws.Range("A101").Value = 'nem'
$ws.Range("B101").Value = '2021-01-11 01:33:28.969'
$ws.Range("A102").Value = 'nttyy'
$ws.Range("B102").Value = '2021-01-11 01:33:28.981'
$ws.Range("A103").Value = 'nsc'
$ws.Range("B103").Value = '2021-01-11 01:33:28.995'
$ws.Range("A104").Value = 'ntr'
$ws.Range("B104").Value = '2021-01-11 01:33:29.007'
$ws.Range("A105").Value = 'onb'
$ws.Range("B105").Value = '2021-01-11 01:33:29.018'
$ws.Range("A106").Value = 'orcl'
$ws.Range("B106").Value = '2021-01-11 01:33:29.029'
$ws.Range("A107").Value = 'pfe'
$ws.Range("B107").Value = '2021-01-11 01:33:29.041'
$ws.Range("A108").Value = 'doc'
$ws.Range("B108").Value = '2021-01-11 01:33:29.052'
$ws.Range("A109").Value = 'pnc'
$ws.Range("B109").Value = '2021-01-11 01:33:29.064'
$ws.Range("A110").Value = 'pru'
$ws.Range("B110").Value = '2021-01-11 01:33:29.074'
$ws.Range("A111").Value = 'qcom'
$ws.Range("B111").Value = '2021-01-11 01:33:29.085'
$ws.Range("A112").Value = 'reg'
$ws.Range("B112").Value = '2021-01-11 01:33:29.097'
$ws.Range("A113").Value = 'rcl'
$ws.Range("B113").Value = '2021-01-11 01:33:29.106'
$ws.Range("A114").Value = 'rds.a'
$ws.Range("B114").Value = '2021-01-11 01:33:29.115'
$ws.Range("A115").Value = 'sny'
$ws.Range("B115").Value = '2021-01-11 01:33:29.125'
$ws.Range("A116").Value = 'slb'
$ws.Range("B116").Value = '2021-01-11 01:33:29.135'
$ws.Range("A117").Value = 'stx'
$ws.Range("B117").Value = '2021-01-11 01:33:29.147'
$ws.Range("A118").Value = 'sfl'
$ws.Range("B118").Value = '2021-01-11 01:33:29.16'
$ws.Range("A119").Value = 'shpg'
$ws.Range("B119").Value = '2021-01-11 01:33:29.172'
$ws.Range("A120").Value = 'syf'
$ws.Range("B120").Value = '2021-01-11 01:33:29.186'
$ws.Range("A121").Value = 'tak'
$ws.Range("B121").Value = '2021-01-11 01:33:29.198'
$ws.Range("A122").Value = 'tgt'
$ws.Range("B122").Value = '2021-01-11 01:33:29.214'
$ws.Range("A123").Value = 'tot'
$ws.Range("B123").Value = '2021-01-11 01:33:29.224'
$ws.Range("A124").Value = 'trn'
$ws.Range("B124").Value = '2021-01-11 01:33:29.236'
$ws.Range("A125").Value = 'tnp'
$ws.Range("B125").Value = '2021-01-11 01:33:29.246'
$ws.Range("A126").Value = 'tpc'
$ws.Range("B126").Value = '2021-01-11 01:33:29.26'
$ws.Range("A127").Value = 'tsn'
$ws.Range("B127").Value = '2021-01-11 01:33:29.27'
$ws.Range("A128").Value = 'wmt'
$ws.Range("B128").Value = '2021-01-11 01:33:29.28'
$ws.Range("A129").Value = 'dis'
$ws.Range("B129").Value = '2021-01-11 01:33:29.289'
$ws.Range("A130").Value = 'wfc'
$ws.Range("B130").Value = '2021-01-11 01:33:29.299'
$ws.Range("A131").Value = 'whr'
$ws.Range("B131").Value = '2021-01-11 01:33:29.308'
$ws.Range("A132").Value = 'wsm'
$ws.Range("B132").Value = '2021-01-11 01:33:29.316'
$ws.Range("A133").Value = 'auy'
$ws.Range("B133").Value = '2021-01-11 01:33:29.327'
$ws.Range("A134").Value = 'zbh'
$ws.Range("B134").Value = '2021-01-11 01:33:29.337'
